$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(196).Insert()

$ws.Cells.Item(196, 1).Value = 7
$ws.Cells.Item(196, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(196, 3).Value = "Ñuble"
$ws.Cells.Item(196, 4).Value = 44628
$ws.Cells.Item(196, 5).Value = 16
$ws.Cells.Item(196, 6).Value = 100112032
$ws.Cells.Item(196, 7).Value = "Zapallo italiano"
$ws.Cells.Item(196, 8).Value = "Sin especificar"
$ws.Cells.Item(196, 9).Value = "Primera"
$ws.Cells.Item(196, 10).Value = 120
$ws.Cells.Item(196, 11).Value = 7500
$ws.Cells.Item(196, 12).Value = 8000
$ws.Cells.Item(196, 13).Value = 7750
$ws.Cells.Item(196, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(196, 15).Value = "Región del Maule"
$ws.Cells.Item(196, 16).Value = 155
$ws.Cells.Item(196, 17).Value = 50
$ws.Cells.Item(196, 18).Value = "Hortaliza"
